$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 809/810, shifting all existing data (rows 809-852)
# down to rows 811-854.
$ws.Rows("809:810").Insert()

# --- New row 809 ---
$ws.Range("A809").Value = 10
$ws.Range("B809").Value = "Vega Modelo de Temuco"
$ws.Range("C809").Value = "La Araucanía"
$ws.Range("D809").Value2 = 45041
$ws.Range("E809").Value = 9
$ws.Range("F809").Value = "Fruta"
$ws.Range("G809").Value = 100102
$ws.Range("H809").Value = "Cítricos"
$ws.Range("I809").Value = 100102004
$ws.Range("J809").Value = "Mandarina"
$ws.Range("K809").Value = "Clementina"
$ws.Range("L809").Value = "Especial"
$ws.Range("M809").Value = 65
$ws.Range("N809").Value = 26000
$ws.Range("O809").Value = 26000
$ws.Range("P809").Value = 26000
$ws.Range("Q809").Value = "`$/caja 15 kilos"
$ws.Range("R809").Value = "Región de O'Higgins"
$ws.Range("S809").Value = 1733
$ws.Range("T809").Value = 15

# --- New row 810 ---
$ws.Range("A810").Value = 10
$ws.Range("B810").Value = "Vega Modelo de Temuco"
$ws.Range("C810").Value = "La Araucanía"
$ws.Range("D810").Value2 = 45041
$ws.Range("E810").Value = 9
$ws.Range("F810").Value = "Fruta"
$ws.Range("G810").Value = 100102
$ws.Range("H810").Value = "Cítricos"
$ws.Range("I810").Value = 100102004
$ws.Range("J810").Value = "Mandarina"
$ws.Range("K810").Value = "Clementina"
$ws.Range("L810").Value = "Primera"
$ws.Range("M810").Value = 155
$ws.Range("N810").Value = 20000
$ws.Range("O810").Value = 20000
$ws.Range("P810").Value = 20000
$ws.Range("Q810").Value = "`$/bandeja 18 kilos"
$ws.Range("R810").Value = "Región de O'Higgins"
$ws.Range("S810").Value = 1111
$ws.Range("T810").Value = 18
